# Update gh-pages to output generated at 456a3b4
# This updates the "want-to-go" counts (column F) for several events across
# the 展览 (Exhibition), 本地生活 (Local life) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3308
$ws1.Range("F12").Value = 26
$ws1.Range("F24").Value = 5881
$ws1.Range("F32").Value = 125
$ws1.Range("F35").Value = 3
$ws1.Range("F36").Value = 780
$ws1.Range("F37").Value = 820

# Sheet: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 218
$ws3.Range("F3").Value = 1089

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 218
$ws4.Range("F4").Value = 1089
$ws4.Range("F8").Value = 3308
$ws4.Range("F16").Value = 26
$ws4.Range("F29").Value = 5881
$ws4.Range("F38").Value = 125
$ws4.Range("F41").Value = 3
$ws4.Range("F42").Value = 780
$ws4.Range("F43").Value = 820
